# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (cloned from "2021-Q4" so it keeps the
#    same header layout / styles) right after "2021-Q4" and before "总计".
# 2. Fill it with the 2022-Q1 per-fund holding data.
# 3. Prepend a "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone the "2021-Q4" sheet to get identical formatting
# (bold/boxed header row, boxed index column) and place it right after it.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)

$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# The source sheet had 4 data rows (rows 2-5); 2022-Q1 only needs 2, so
# drop the extra two rows (this also fixes up the <dimension> to A1:H3).
$newSheet.Rows("4:5").Delete()

# ---------------------------------------------------------------------
# Step 2: write the 2022-Q1 fund table. Fund codes / amounts are numeric
# looking strings that must stay text (leading zeros, fixed decimals), so
# force a text number format before assigning them; the ranking column
# (H) is a genuine number.
# ---------------------------------------------------------------------
$textCols = "B2:G3"
$newSheet.Range($textCols).NumberFormat = "@"

$newSheet.Range("B2").Value = "004738"
$newSheet.Range("C2").Value = "上投摩根安隆回报混合A"
$newSheet.Range("D2").Value = "23.04"
$newSheet.Range("E2").Value = "21.18"
$newSheet.Range("F2").Value = "0.80"
$newSheet.Range("G2").Value = "0.1843"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "004739"
$newSheet.Range("C3").Value = "上投摩根安隆回报混合C"
$newSheet.Range("D3").Value = "7.32"
$newSheet.Range("E3").Value = "21.18"
$newSheet.Range("F3").Value = "0.80"
$newSheet.Range("G3").Value = "0.0586"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# Step 3: update "总计" with a new first data row for 2022-Q1, pushing
# the existing rows down by one.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Rows("2:2").Insert()
$zj.Range("B2:D2").ClearFormats()

# Re-apply the boxed index-column style to A2 by copying (formats only)
# from A3 below - it already carries the exact style used throughout the
# sheet, so this reuses it instead of registering a new one.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.24

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
$zj.Range("A6").Value = 4
